# Applies the three changes described by the diff:
#   1. "Ticket system for tracking repairings needs and monitor the progress"
#      -> wrap "repairings" in <w:proofErr> spellStart/spellEnd (split run).
#   2. "Mouisture sensors " paragraph
#      -> wrap "Mouisture" in <w:proofErr> spellStart/spellEnd (split run).
#   3. Delete the whole "Machine learning algorithms..." bullet paragraph and
#      trim "GUI for desktop and mobile" down to "GUI for desktop ".

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Change 1: "Ticket " + "system for tracking repairings" + " needs..."
#           -> "Ticket " + "system for tracking " + proofErr(repairings) + " needs..."
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ticket system for tracking repairings needs and monitor the progress*") {
        $body = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
                '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Ticket </w:t></w:r>' + `
                '<w:r w:rsidR="006C44AE"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">system for tracking </w:t></w:r>' + `
                '<w:proofErr w:type="spellStart"/>' + `
                '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>repairings</w:t></w:r>' + `
                '<w:proofErr w:type="spellEnd"/>' + `
                '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> needs and monitor the progress</w:t></w:r>' + `
                '</w:p>'
        $p.Range.InsertXML($pkgOpen + $body + $pkgClose)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 2: "Mouisture " + "sensors" + " "
#           -> proofErr(Mouisture) + " " + "sensors" + " "
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Mouisture sensors*") {
        $body = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
                '<w:proofErr w:type="spellStart"/>' + `
                '<w:r w:rsidRPr="009516DE"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Mouisture</w:t></w:r>' + `
                '<w:proofErr w:type="spellEnd"/>' + `
                '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
                '<w:r w:rsidR="006C44AE"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>sensors</w:t></w:r>' + `
                '<w:r w:rsidRPr="009516DE"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
                '</w:p>'
        $p.Range.InsertXML($pkgOpen + $body + $pkgClose)
        break
    }
}

# ---------------------------------------------------------------------------
# Change 3a: delete the whole "Machine learning algorithms..." paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Machine learning algorithms to detect irregularities in the measured data*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Change 3b: "GUI for desktop and mobile" -> "GUI for desktop "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("GUI for desktop and mobile", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "GUI for desktop ", 2) | Out-Null
